$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# Header cell F1: S000 -> S020
$ws.Range("F1").Value = "S020"

# Row 3
$ws.Range("C3").Value = "no walk/same"
$ws.Range("E3").Value = "start DD"
$ws.Range("I3").Value = "walk/same"
$ws.Range("K3").Value = "start SD"

# Row 10
$ws.Range("C10").Value = "walk/diff"
$ws.Range("E10").Value = "start SD"
$ws.Range("I10").Value = "walk/diff"
$ws.Range("K10").Value = "start DD"

# Row 17
$ws.Range("C17").Value = "no walk/diff"
$ws.Range("E17").Value = "start SD"
$ws.Range("I17").Value = "no walk/diff"
$ws.Range("K17").Value = "start SD"

# Row 24
$ws.Range("C24").Value = "walk/same"
$ws.Range("E24").Value = "start DD"
$ws.Range("I24").Value = "walk/diff"
$ws.Range("K24").Value = "start SD"

# Row 31
$ws.Range("C31").Value = "no walk/same"
$ws.Range("E31").Value = "start DD"
$ws.Range("I31").Value = "no walk/same"
$ws.Range("K31").Value = "start DD"
